$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "5533988797979"
$ws.Range("B3").Value = "5533988797979"
$ws.Range("B4").Value = "5533988797979"
$ws.Range("B5").Value = "5533988797979"
$ws.Range("B6").Value = "5533988797979"

$ws.Columns.Item(2).ColumnWidth = 29.7109375

$ws.Range("B5").Select()
